$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "43.751.63"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -0.46%  "
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.348.53"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").Value = "  -0.30%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "233.84"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("E6").Value = "  +2.34%  "
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "66.10"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  +4.34%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +1.48%  "
$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0973"
$cell.Style = $origStyle
$ws.Range("E10").Value = "  -4.15%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "56.66"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -0.50%  "
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "26.89"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  +2.27%  "
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.697.79"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +3.49%  "
$ws.Range("E14").Value = "  -0.95%  "
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.46"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("E16").Value = "  -0.42%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.858"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +1.92%  "
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.354.27"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  +3.97%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "43.694.13"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -0.47%  "
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0983"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -2.19%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "74.07"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("E22").Value = "  +3.23%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "249.41"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("B24").Value = "WEMIXToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.81"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  +14.79%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  -0.05%  "
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("E28").Value = "  -0.74%  "
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.41"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  +7.67%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "174.88"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  +3.87%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.130"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -5.93%  "
$ws.Range("E34").Value = "  +4.49%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0688"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  -1.30%  "
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.02"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +2.38%  "
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.72"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  -2.53%  "
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  +5.84%  "
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.57"
$cell.Style = $origStyle
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("E40").Value = "  -2.14%  "
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.21"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  +11.68%  "
$ws.Range("E42").Value = "  -0.13%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "17.98"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("E44").Value = "  +10.59%  "
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "99.64"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  +1.74%  "
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0954"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -2.06%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.450.59"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.98"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  -1.48%  "
$ws.Range("E51").Value = "  +0.49%  "
